$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full roster (header stays the same; data rows are re-sorted with one
# brand-new player "Spencer Dinwiddie" inserted at the top).
$data = @(
    @("Spencer Dinwiddie",        "PG,SG",      "Dallas Mavericks"),
    @("Jordan Poole",             "PG,SG",      "Washington Wizards"),
    @("CJ McCollum",              "PG,SG",      "New Orleans Pelicans"),
    @("Klay Thompson",            "SG,SF",      "Dallas Mavericks"),
    @("Lauri Markkanen",          "SF,PF",      "Utah Jazz"),
    @("Quentin Grimes",           "SG,SF",      "Dallas Mavericks"),
    @("Keyonte George",           "PG,SG",      "Utah Jazz"),
    @("John Collins",             "PF,C",       "Utah Jazz"),
    @("Tobias Harris",            "SF,PF",      "Detroit Pistons"),
    @("Zach LaVine",              "SG,SF",      "Chicago Bulls"),
    @("Joel Embiid",              "C",          "Philadelphia 76ers"),
    @("Guerschon Yabusele",       "PF,C",       "Philadelphia 76ers"),
    @("Andrew Nembhard",          "PG,SG",      "Indiana Pacers"),
    @("Shai Gilgeous-Alexander",  "PG,SG",      "Oklahoma City Thunder"),
    @("Jalen Williams",           "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Kyrie Irving",             "PG,SG",      "Dallas Mavericks"),
    @("RJ Barrett",               "SG,SF,PF",   "Toronto Raptors"),
    @("Jimmy Butler",             "SF,PF",      "Miami Heat")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
